# Update cryptos list values (prices and volume %) to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.535.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.227.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "268.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.25%  "
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0921"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +17.22%  "
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.564.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.230.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.797"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.512.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000103"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.59%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.49%  "
$ws.Range("E27").Value = "  +10.56%  "
$ws.Range("E28").Value = "  +5.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0924"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0349"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +20.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.218"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.43%  "
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.434"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.453.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.35%  "
